$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.984.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.643.04'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.54%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4765'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2597'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06099'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07028'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.649.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5886'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.90%  '
$ws.Range('E14').Value = '  -7.37%  '
$ws.Range('E15').Value = '  -5.15%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9997'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '24.991.34'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006592'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('E20').Value = '  -6.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.858.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.287'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.557'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.239'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.04'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.383'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.39%  '
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.632'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.891'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07652'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.584'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9997'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04286'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.574'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9257'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5888'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.580'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8711'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.93%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01504'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.756'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3696'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.680'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1100'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.31%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.092'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.64%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05205'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '28.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9994'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.10%  '
